$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = [double]"2"
$ws.Range("F2").Value = [double]"0.6666666666666666"
$ws.Range("G2").Value = [double]"0.01848533333333334"
$ws.Range("H2").Value = [double]"0.05545600000000001"
$ws.Range("I2").Value = [double]"0.001625201930372746"
$ws.Range("J2").Value = [double]"0.001625201930372746"
$ws.Range("K2").Value = [double]"1"
$ws.Range("L2").Value = [double]"0.6666666666666666"
$ws.Range("M2").Value = [double]"0.05354133333333333"
$ws.Range("N2").Value = [double]"0.160624"
$ws.Range("O2").Value = [double]"0.00209946492164722"
$ws.Range("P2").Value = [double]"0.00209946492164722"
$ws.Range("Q2").Value = [double]"0.0009897293937777779"
$ws.Range("R2").Value = [double]"0.008907564544000001"
$ws.Range("S2").Value = [double]"3.412054443410927E-06"
$ws.Range("T2").Value = [double]"3.412054443410927E-06"

# Row 3
$ws.Range("E3").Value = [double]"2"
$ws.Range("F3").Value = [double]"0.6666666666666666"
$ws.Range("G3").Value = [double]"0.01848533333333334"
$ws.Range("H3").Value = [double]"0.05545600000000001"
$ws.Range("I3").Value = [double]"0.001625201930372746"
$ws.Range("J3").Value = [double]"0.001625201930372746"
$ws.Range("K3").Value = [double]"3"
$ws.Range("L3").Value = [double]"1"
$ws.Range("M3").Value = [double]"1.392600333333333"
$ws.Range("N3").Value = [double]"4.177801000000001"
$ws.Range("O3").Value = [double]"0.05460670042535784"
$ws.Range("P3").Value = [double]"0.05460670042535784"
$ws.Range("Q3").Value = [double]"0.02574268136177778"
$ws.Range("R3").Value = [double]"0.231684132256"
$ws.Range("S3").Value = [double]"8.87469149425778E-05"
$ws.Range("T3").Value = [double]"8.874691494257779E-05"

# Row 4
$ws.Range("E4").Value = [double]"2"
$ws.Range("F4").Value = [double]"0.6666666666666666"
$ws.Range("G4").Value = [double]"0.01848533333333334"
$ws.Range("H4").Value = [double]"0.05545600000000001"
$ws.Range("I4").Value = [double]"0.001625201930372746"
$ws.Range("J4").Value = [double]"0.001625201930372746"
$ws.Range("K4").Value = [double]"3"
$ws.Range("L4").Value = [double]"1"
$ws.Range("M4").Value = [double]"24.05622933333333"
$ws.Range("N4").Value = [double]"72.168688"
$ws.Range("O4").Value = [double]"0.943293834652995"
$ws.Range("P4").Value = [double]"0.943293834652995"
$ws.Range("Q4").Value = [double]"0.4446874179697778"
$ws.Range("R4").Value = [double]"4.002186761728001"
$ws.Range("S4").Value = [double]"0.001533042960986757"
$ws.Range("T4").Value = [double]"0.001533042960986757"

# Row 5
$ws.Range("E5").Value = [double]"2"
$ws.Range("F5").Value = [double]"0.6666666666666666"
$ws.Range("G5").Value = [double]"0.03069133333333333"
$ws.Range("H5").Value = [double]"0.092074"
$ws.Range("I5").Value = [double]"0.002698334581238102"
$ws.Range("J5").Value = [double]"0.002698334581238102"
$ws.Range("K5").Value = [double]"1"
$ws.Range("L5").Value = [double]"0.3333333333333333"
$ws.Range("M5").Value = [double]"0.05354133333333333"
$ws.Range("N5").Value = [double]"0.160624"
$ws.Range("O5").Value = [double]"0.00209946492164722"
$ws.Range("P5").Value = [double]"0.00209946492164722"
$ws.Range("Q5").Value = [double]"0.001643254908444444"
$ws.Range("R5").Value = [double]"0.014789294176"
$ws.Range("S5").Value = [double]"5.665058800177035E-06"
$ws.Range("T5").Value = [double]"5.665058800177034E-06"

# Row 6
$ws.Range("E6").Value = [double]"2"
$ws.Range("F6").Value = [double]"0.6666666666666666"
$ws.Range("G6").Value = [double]"0.03069133333333333"
$ws.Range("H6").Value = [double]"0.092074"
$ws.Range("I6").Value = [double]"0.002698334581238102"
$ws.Range("J6").Value = [double]"0.002698334581238102"
$ws.Range("K6").Value = [double]"3"
$ws.Range("L6").Value = [double]"1"
$ws.Range("M6").Value = [double]"1.392600333333333"
$ws.Range("N6").Value = [double]"4.177801000000001"
$ws.Range("O6").Value = [double]"0.05460670042535784"
$ws.Range("P6").Value = [double]"0.05460670042535784"
$ws.Range("Q6").Value = [double]"0.04274076103044445"
$ws.Range("R6").Value = [double]"0.3846668492740001"
$ws.Range("S6").Value = [double]"0.0001473471481250524"
$ws.Range("T6").Value = [double]"0.0001473471481250524"

# Row 7
$ws.Range("E7").Value = [double]"2"
$ws.Range("F7").Value = [double]"0.6666666666666666"
$ws.Range("G7").Value = [double]"0.03069133333333333"
$ws.Range("H7").Value = [double]"0.092074"
$ws.Range("I7").Value = [double]"0.002698334581238102"
$ws.Range("J7").Value = [double]"0.002698334581238102"
$ws.Range("K7").Value = [double]"3"
$ws.Range("L7").Value = [double]"1"
$ws.Range("M7").Value = [double]"24.05622933333333"
$ws.Range("N7").Value = [double]"72.168688"
$ws.Range("O7").Value = [double]"0.943293834652995"
$ws.Range("P7").Value = [double]"0.943293834652995"
$ws.Range("Q7").Value = [double]"0.7383177532124445"
$ws.Range("R7").Value = [double]"6.644859778912"
$ws.Range("S7").Value = [double]"0.002545322374312873"
$ws.Range("T7").Value = [double]"0.002545322374312872"

# Row 8
$ws.Range("E8").Value = [double]"3"
$ws.Range("F8").Value = [double]"1"
$ws.Range("G8").Value = [double]"11.32499966666667"
$ws.Range("H8").Value = [double]"33.974999"
$ws.Range("I8").Value = [double]"0.9956764634883892"
$ws.Range("J8").Value = [double]"0.995676463488389"
$ws.Range("K8").Value = [double]"1"
$ws.Range("L8").Value = [double]"0.3333333333333333"
$ws.Range("M8").Value = [double]"0.05354133333333333"
$ws.Range("N8").Value = [double]"0.160624"
$ws.Range("O8").Value = [double]"0.00209946492164722"
$ws.Range("P8").Value = [double]"0.00209946492164722"
$ws.Range("Q8").Value = [double]"0.606355582152889"
$ws.Range("R8").Value = [double]"5.457200239376"
$ws.Range("S8").Value = [double]"0.002090387808403632"
$ws.Range("T8").Value = [double]"0.002090387808403631"

# Row 9
$ws.Range("E9").Value = [double]"3"
$ws.Range("F9").Value = [double]"1"
$ws.Range("G9").Value = [double]"11.32499966666667"
$ws.Range("H9").Value = [double]"33.974999"
$ws.Range("I9").Value = [double]"0.9956764634883892"
$ws.Range("J9").Value = [double]"0.995676463488389"
$ws.Range("K9").Value = [double]"3"
$ws.Range("L9").Value = [double]"1"
$ws.Range("M9").Value = [double]"1.392600333333333"
$ws.Range("N9").Value = [double]"4.177801000000001"
$ws.Range("O9").Value = [double]"0.05460670042535784"
$ws.Range("P9").Value = [double]"0.05460670042535784"
$ws.Range("Q9").Value = [double]"15.77119831079989"
$ws.Range("R9").Value = [double]"141.940784797199"
$ws.Range("S9").Value = [double]"0.05437060636229021"
$ws.Range("T9").Value = [double]"0.0543706063622902"

# Row 10
$ws.Range("E10").Value = [double]"3"
$ws.Range("F10").Value = [double]"1"
$ws.Range("G10").Value = [double]"11.32499966666667"
$ws.Range("H10").Value = [double]"33.974999"
$ws.Range("I10").Value = [double]"0.9956764634883892"
$ws.Range("J10").Value = [double]"0.995676463488389"
$ws.Range("K10").Value = [double]"3"
$ws.Range("L10").Value = [double]"1"
$ws.Range("M10").Value = [double]"24.05622933333333"
$ws.Range("N10").Value = [double]"72.168688"
$ws.Range("O10").Value = [double]"0.943293834652995"
$ws.Range("P10").Value = [double]"0.943293834652995"
$ws.Range("Q10").Value = [double]"272.4367891812569"
$ws.Range("R10").Value = [double]"2451.931102631313"
$ws.Range("S10").Value = [double]"0.9392154693176954"
$ws.Range("T10").Value = [double]"0.9392154693176952"
